$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Body text edit: add two new paragraphs ("Husd" and "jksds")
#    right before the trailing bookmark ("_GoBack") paragraph, i.e.
#    right after the "Ravan" paragraph.
# ------------------------------------------------------------------

# Locate the "Ravan" paragraph (last paragraph with visible text,
# right before the empty bookmark-only paragraph).
$ravanPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.Trim() -eq "Ravan") {
        $ravanPara = $p
    }
}

# Insert a brand new paragraph right after "Ravan" and give it the
# text "Husd".
$ravanPara.Range.InsertParagraphAfter()
$husdPara = $ravanPara.Next()
$husdPara.Range.Text = "Husd"

# The paragraph that used to directly follow "Ravan" now follows
# "Husd" -- it's the one holding the _GoBack bookmark. Insert the
# "jksds" run in front of the bookmark, inside that same paragraph.
$bookmarkPara = $husdPara.Next()
$bookmarkPara.Range.InsertBefore("jksds")

# ------------------------------------------------------------------
# 2) Style edit: mark the built-in "Normal Table" style as a quick
#    style (w:qFormat) -- equivalent of checking "Add to Quick Style
#    list" in the Modify Style dialog.
# ------------------------------------------------------------------

$tableStyle = $d.Styles.Item("Normal Table")
$tableStyle.QuickStyle = $true
